$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.732.32"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.091.01"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.96"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.38"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.620.56"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.73"
$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.744.95"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.087.28"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.14"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  -1.49%  "

$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "344.87"
$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.505"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.73"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0922"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("E28").Value = "  +2.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.01"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("E32").Value = "  +1.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.42"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.61"
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.12"
$ws.Range("E35").Value = "  +3.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.78"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").Value = "  +3.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0685"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.130.05"
$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.74"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("E44").Value = "  +3.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.286.10"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.961"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.99"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("E50").Value = "  +7.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "264.14"
$ws.Range("E51").Value = "  +11.41%  "
